# Weekly update for "Fruta, Terminal La Palmera de La Serena - Coco"
# - updates D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
#   P (Precio promedio ponderado), S (Precio $/Kg) for existing rows 5-13
# - appends four new data rows (14-17) with the same row layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("D$Row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMin
    $ws.Range("O$Row").Value = $PrecioMax
    $ws.Range("P$Row").Value = $PrecioProm
    $ws.Range("S$Row").Value = $PrecioKg
}

function New-DataRow {
    param($Row, $Fecha, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Range("A$Row").Value = 8
    $ws.Range("B$Row").Value = "Terminal La Palmera de La Serena"
    $ws.Range("C$Row").Value = "Coquimbo"
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("D$Row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("E$Row").Value = 4
    $ws.Range("F$Row").Value = "Fruta"
    $ws.Range("G$Row").Value = 100108
    $ws.Range("H$Row").Value = "Tropicales y subtropicales"
    $ws.Range("I$Row").Value = 100108007
    $ws.Range("J$Row").Value = "Coco"
    $ws.Range("K$Row").Value = "Sin especificar"
    $ws.Range("L$Row").Value = "Primera"
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMin
    $ws.Range("O$Row").Value = $PrecioMax
    $ws.Range("P$Row").Value = $PrecioProm
    $ws.Range("Q$Row").Value = "`$/malla 20 unidades"
    $ws.Range("R$Row").Value = "Perú"
    $ws.Range("S$Row").Value = $PrecioKg
    $ws.Range("T$Row").Value = 20
}

# --- updates to existing rows ---
$ws.Range("D5").Value = 44427
$ws.Range("D5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Set-DataRow 6  44410 200 20000 21000 20500 1025
Set-DataRow 7  44336 100 19500 20000 19750 988
Set-DataRow 8  44364 140 20000 21000 20500 1025
Set-DataRow 9  44350 160 19000 20000 19500 975
Set-DataRow 10 44315 100 20000 21000 20500 1025
Set-DataRow 11 44417 160 20000 21000 20500 1025
Set-DataRow 12 44420 160 20000 21000 20500 1025

$ws.Range("D13").Value = 44333
$ws.Range("D13").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- newly appended rows ---
New-DataRow 14 44301 100 18000 19000 18500 925
New-DataRow 15 44326 160 19500 20000 19750 988
New-DataRow 16 44343 100 19500 20000 19750 988
New-DataRow 17 44418 200 20000 21000 20500 1025
